# Applies the "Issue 7 2 (#15)" changes to shareyourcloning_linkml.xlsx:
#  * RepositoryIdSource / BenchlingUrlSource: swap repository_name/repository_id
#    column order (id first, name second) and move their dropdown validation
#    from column A to column B accordingly.
#  * AddGeneIdSource: same repository_id/repository_name swap, but in columns
#    C/D, moving the dropdown validation from column C to column D.
#  * AssemblyJoinComponent -> renamed to AssemblyFragment, "location" column
#    split into nullable "left_location" and "right_location" columns.
#  * AssemblyJoin sheet removed entirely (no longer part of the schema).

$wb = $excel.ActiveWorkbook

# --- RepositoryIdSource: repository_name/repository_id -> repository_id/repository_name ---
$ws = $wb.Worksheets.Item("RepositoryIdSource")
$ws.Range("A1").Value = "repository_id"
$ws.Range("B1").Value = "repository_name"
$ws.Range("A2:A1048576").Validation.Delete()
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')

# --- AddGeneIdSource: repository_name/repository_id -> repository_id/repository_name (cols C/D) ---
$ws = $wb.Worksheets.Item("AddGeneIdSource")
$ws.Range("C1").Value = "repository_id"
$ws.Range("D1").Value = "repository_name"
$ws.Range("C2:C1048576").Validation.Delete()
$ws.Range("D2:D1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')

# --- BenchlingUrlSource: repository_name/repository_id -> repository_id/repository_name ---
$ws = $wb.Worksheets.Item("BenchlingUrlSource")
$ws.Range("A1").Value = "repository_id"
$ws.Range("B1").Value = "repository_name"
$ws.Range("A2:A1048576").Validation.Delete()
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')

# --- AssemblyJoinComponent -> AssemblyFragment: location -> left_location + right_location ---
$ws = $wb.Worksheets.Item("AssemblyJoinComponent")
$ws.Name = "AssemblyFragment"
$ws.Range("B1").Value = "left_location"
$ws.Range("C1").Value = "right_location"
$ws.Range("D1").Value = "reverse_complemented"

# --- AssemblyJoin sheet is no longer part of the schema: remove it ---
$wb.Worksheets.Item("AssemblyJoin").Delete() | Out-Null
